$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.225.66'
$ws.Range("E2").Value = '  +0.84%  '
$ws.Range("D3").Value = '2.933.83'
$ws.Range("E3").Value = '  +4.81%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '''354.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.03%  '
$ws.Range("D6").Value = '''113.42'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").Value = '''0.560'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '''0.625'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.01%  '
$ws.Range("D10").Value = '''39.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("D11").Value = '''0.0884'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.43%  '
$ws.Range("D12").Value = '''0.137'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("D13").Value = '''20.09'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("D14").Value = '''7.79'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.19%  '
$ws.Range("D15").Value = '3.374.16'
$ws.Range("E15").Value = '  +4.22%  '
$ws.Range("D16").Value = '2.901.69'
$ws.Range("E16").Value = '  +4.07%  '
$ws.Range("D17").Value = '''0.990'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.45%  '
$ws.Range("D18").Value = '52.241.08'
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("E19").Value = '  -1.97%  '
$ws.Range("D20").Value = '''7.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.13%  '
$ws.Range("D21").Value = '''14.19'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.70%  '
$ws.Range("D22").Value = '0.0₃0983'
$ws.Range("E22").Value = '  +1.00%  '
$ws.Range("D23").Value = '''71.30'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.26%  '
$ws.Range("D24").Value = '''270.69'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").Value = '''2.82'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.42%  '
$ws.Range("E26").Value = '  +11.66%  '
$ws.Range("D27").Value = '''26.98'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.25%  '
$ws.Range("E28").Value = '  +0.04%  '
$ws.Range("D29").Value = '''7.05'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +15.42%  '
$ws.Range("D30").Value = '''10.68'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.33%  '
$ws.Range("E31").Value = '  +14.71%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("D33").Value = '''37.32'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.97%  '
$ws.Range("D34").Value = '''6.04'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.17%  '
$ws.Range("D35").Value = '''53.19'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.81%  '
$ws.Range("D36").Value = '''0.0454'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.21%  '
$ws.Range("D37").Value = '''0.998'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '''3.38'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.69%  '
$ws.Range("D39").Value = '''18.83'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '''2.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.44%  '
$ws.Range("D41").Value = '''2.72'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.92%  '
$ws.Range("D42").Value = '''0.118'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.15%  '
$ws.Range("D43").Value = '''23.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.94%  '
$ws.Range("D44").Value = '''2.19'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.63%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '''3.55'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D46").Value = '''2.53'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.37%  '
$ws.Range("B47").Value = 'Maker'
$ws.Range("C47").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D47").Value = '2.192.38'
$ws.Range("E47").Value = '  +3.33%  '
$ws.Range("B48").Value = 'Monero'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D48").Value = '''116.27'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.70%  '
$ws.Range("D49").Value = '''0.252'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +13.66%  '
$ws.Range("E50").Value = '  +10.34%  '
$ws.Range("D51").Value = '''0.959'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.50%  '
